# Apply the edit described by the commit diff:
#  1. The writer placed the cursor inside the word "influenced" (right after
#     "Kishida was i") and kept typing/re-typing there, which is why Word's
#     "_GoBack" last-edit-position bookmark moved from its old spot (right
#     after the opening "Kishida,") to that exact point, splitting the run
#     that used to hold the whole sentence into two runs.
#  2. The header gained an affiliation after the author's name.

$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark into the middle of the sentence -------
# Find the split point: right after "Kishida was i" / right before "nfluenced".
$rng = $d.Content
$null = $rng.Find.Execute("Kishida was i", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Adding a bookmark named "_GoBack" automatically replaces any existing one
# of that name (Word only ever keeps a single "_GoBack" bookmark), so this
# both removes the old bookmark (next to "Kishida,") and creates the new one
# here -- which also splits the run in two, matching the diff.
$d.Bookmarks.Add("_GoBack", $rng)

# --- 2. Add the affiliation text in the header ----------------------------
$hdr = $d.Sections.First.Headers.Item(1)
$hrng = $hdr.Range
$hrng.Collapse(0)
$hrng.InsertAfter(", National University of Singapore")
